# Apply the edit described by the diff: swap the species/taxon-related
# values between row 10 and row 11 (columns A, B, E, F, G, H, Q, R),
# leaving all other columns (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB,
# AD, AE, AG, AT, AW, AX, AY) untouched since they already match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 -> new values (previously held by row 11)
$ws.Range("A10").Value = 111519524
$ws.Range("B10").Value = 77515
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 633583.7615760232
$ws.Range("R10").Value = 7117850.915647855

# Row 11 -> new values (previously held by row 10)
$ws.Range("A11").Value = 111519523
$ws.Range("B11").Value = 89686
$ws.Range("E11").Value = 658
$ws.Range("F11").Value = "Rosenticka"
$ws.Range("G11").Value = "Rhodofomes roseus"
$ws.Range("H11").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q11").Value = 633714.5983269843
$ws.Range("R11").Value = 7117626.805168894
